# Add the new "ESP 32 S3 Wroom N16R8" BOM line (row 9) to the BOM sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the new row's cells. Order matches the shared-string insertion order
# in the target workbook: Comment, Footprint, Designator, JLCPCB Part#.
$ws.Range("A9").Value = "ESP 32 S3 Wroom N16R8"
$ws.Range("C9").Value = "ESP32-S3-WROOM-1-N16R8"
$ws.Range("B9").Value = "X1"
$ws.Range("D9").Value = "C2913202"

# Resize column A to fit the new (longer) Comment text.
$ws.Columns.Item(1).ColumnWidth = 20.166666

# Match the selection left behind after entering the new data.
$ws.Range("D9").Select() | Out-Null

# Set page orientation (as recorded in the saved page setup).
$ws.PageSetup.Orientation = 1
